$d = $word.ActiveDocument
foreach ($p in $d.Paragraphs) {
    $p.Range.ParagraphFormat.ContextualSpacing = $false
}
